# Update "Work Week and Social Spending" GDP per Capita data for Algeria
# (Country Code 12, sheet "Data") and append new years 2011-2016.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- 1. Update existing "Data" values (rows 2-192, column E) ---
# These are stored as text (numeric-looking strings), so a leading
# apostrophe is used to force text entry and avoid them being
# reinterpreted as numbers.
$rowsToUpdate = @(2, 52, 95, 132, 133, 134, 135, 136, 137, 138, 139, 140, 141, 142, 143, 144, 145, 146, 147, 148, 149, 150, 151, 152, 153, 154, 155, 156, 157, 158, 159, 160, 161, 162, 163, 164, 165, 166, 167, 168, 169, 170, 171, 172, 173, 174, 175, 176, 177, 178, 179, 180, 181, 182, 183, 184, 185, 186, 187, 188, 189, 190, 191, 192)
$newValues = @("685", "1140", "1854", "2176", "2147", "2193", "2182", "2291", "2303", "2475", "2699", "2740", "3178", "3328", "2868", "2284", "2818", "2879", "2981", "2750", "2907", "3151", "3355", "3585", "3188", "3746", "3757", "3870", "4020", "4157", "4398", "4812", "5088", "5024", "4991", "5145", "5243", "5361", "5469", "5262", "5088", "4850", "4889", "4697", "4708.23191341973", "4850.22333074477", "4816.92359011747", "4856.98538931291", "5142.38400437617", "5451.63444208551", "5638.48016713353", "6069.77731672431", "6421.508468349", "6834.55387976764", "7218.24576822021", "7814.83816030878", "8590.17752627433", "9182.29418389124", "9969.62260937367", "10385.7791083296", "10974.7422866613", "11475.0024058251", "11907.021496637", "12587.744251356")

for ($i = 0; $i -lt $rowsToUpdate.Length; $i++) {
    $r = $rowsToUpdate[$i]
    $v = $newValues[$i]
    $ws.Cells.Item($r, 5).Value = "'" + $v
}

# --- 2. Append new rows for years 2011-2016 ---
$newRowNums = @(193, 194, 195, 196, 197, 198)
$newRowYears = @(2011, 2012, 2013, 2014, 2015, 2016)
$newRowVals = @("13204", "13379", "13494", "13744", "14004", "14331")

for ($i = 0; $i -lt $newRowNums.Length; $i++) {
    $r = $newRowNums[$i]
    $ws.Cells.Item($r, 1).Value = 12
    $ws.Cells.Item($r, 2).Value = "Algeria"
    $ws.Cells.Item($r, 3).Value = "GDP per Capita"
    $ws.Cells.Item($r, 4).Value = $newRowYears[$i]
    $ws.Cells.Item($r, 5).Value = "'" + $newRowVals[$i]
}
